$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 128
$ws.Range("E10").Value = 449
$ws.Range("E23").Value = 184
$ws.Range("E28").Value = 184
$ws.Range("E39").Value = 168
$ws.Range("E42").Value = 331
$ws.Range("F42").Value = 182
$ws.Range("H42").Value = 182
$ws.Range("E46").Value = 285
